# Update BAU production/imports/exports balancing priorities on the FPIEBP sheet.
# Columns: B = production priority, C = imports priority, D = exports priority.
# Most fuel rows move from (production=1, imports=3, exports=2)
# to (production=2, imports=1, exports=3) -- i.e. imports becomes the
# first-priority balancing term, production becomes second, exports last
# (this effectively drives lignite and natural gas exports toward 0, since
# exports is now the lowest-priority / last-resort balancing variable).
# "biofuel diesel" (row 13) moves from (1,2,3) to (2,1,3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPIEBP")

# row -> (production, imports, exports, clearFormatCells)
$updates = @(
    @{ Row = 3;  Vals = @(2,1,3); ClearCols = @() }                 # hard coal
    @{ Row = 4;  Vals = @(2,1,3); ClearCols = @() }                 # natural gas
    @{ Row = 10; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # petroleum gasoline
    @{ Row = 11; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # petroleum diesel
    @{ Row = 12; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # biofuel gasoline
    @{ Row = 13; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # biofuel diesel
    @{ Row = 14; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # jet fuel or kerosene
    @{ Row = 17; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # lignite
    @{ Row = 19; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # heavy fuel oil
    @{ Row = 20; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # LPG propane or butane
    @{ Row = 21; Vals = @(2,1,3); ClearCols = @(2,3,4) }            # municipal solid waste
    @{ Row = 22; Vals = @(2,1,3); ClearCols = @(2) }                # hydrogen (only B22 was styled)
)

foreach ($u in $updates) {
    $row = $u.Row
    $vals = $u.Vals
    foreach ($col in $u.ClearCols) {
        $ws.Cells.Item($row, $col).ClearFormats()
    }
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
